# Update the Shefan (methanol) yearly income-statement dollar workbook:
#  - roll the reporting periods forward by one column (drop the oldest
#    12-month period / disclosure date, add the newest one)
#  - shift all financial figures one column to the left and populate the
#    newest (rightmost) column with the newly reported figures
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 8: "12 ماهه منتهی به ..." period headers (columns D:H) ----
$ws.Range("D8").Value = "12 ماهه منتهی به 1397/12"
$ws.Range("E8").Value = "12 ماهه منتهی به 1398/12"
$ws.Range("F8").Value = "12 ماهه منتهی به 1399/12"
$ws.Range("G8").Value = "12 ماهه منتهی به 1400/12"
$ws.Range("H8").Value = "12 ماهه منتهی به 1401/12"

# ---- Row 9: "تاریخ انتشار" publish-date headers (columns D:H) ----
$ws.Range("D9").Value = "1399-04-04 (8)"
$ws.Range("E9").Value = "1400-04-05 (11)"
$ws.Range("F9").Value = "1401-04-19 (12)"
$ws.Range("G9").Value = "1402-02-30 (9)"
$ws.Range("H9").Value = "1402-02-30"

# ---- Row 11: فروش (Sales) ----
$ws.Range("D11").Value = 262018
$ws.Range("E11").Value = 227202
$ws.Range("F11").Value = 244181
$ws.Range("G11").Value = 356254
$ws.Range("H11").Value = 234907

# ---- Row 12: بهای تمام شده کالای فروش رفته (Cost of goods sold) ----
$ws.Range("D12").Value = -105300
$ws.Range("E12").Value = -102673
$ws.Range("F12").Value = -105352
$ws.Range("G12").Value = -222848
$ws.Range("H12").Value = -171982

# ---- Row 13: سود (زیان) ناخالص (Gross profit) ----
$ws.Range("D13").Value = 156718
$ws.Range("E13").Value = 124529
$ws.Range("F13").Value = 138829
$ws.Range("G13").Value = 133406
$ws.Range("H13").Value = 62924

# ---- Row 14: هزینه های عمومی, اداری و تشکیلاتی (G&A expenses) ----
$ws.Range("D14").Value = -21126
$ws.Range("E14").Value = -37708
$ws.Range("F14").Value = -35363
$ws.Range("G14").Value = -22677
$ws.Range("H14").Value = -26720

# Row 15 (هزینه کاهش ارزش دریافتنی‌ها) stays "-" across the board - no change

# ---- Row 16: خالص سایر درامدها (هزینه ها) ی عملیاتی ----
$ws.Range("D16").Value = 12890
$ws.Range("E16").Value = 22859
$ws.Range("F16").Value = 11559
$ws.Range("G16").Value = 2374
$ws.Range("H16").Value = 369

# ---- Row 17: سود (زیان) عملیاتی (Operating profit) ----
$ws.Range("D17").Value = 148481
$ws.Range("E17").Value = 109679
$ws.Range("F17").Value = 115025
$ws.Range("G17").Value = 113103
$ws.Range("H17").Value = 36574

# ---- Row 18: هزینه های مالی (Financial expenses) ----
$ws.Range("D18").Value = -2350
$ws.Range("E18").Value = -12990
$ws.Range("F18").Value = -11862
$ws.Range("G18").Value = -12260
$ws.Range("H18").Value = -22105

# ---- Row 19: خالص سایر درامدها و هزینه های غیرعملیاتی ----
$ws.Range("D19").Value = 12508
$ws.Range("E19").Value = 28466
$ws.Range("F19").Value = 167404
$ws.Range("G19").Value = -6319
$ws.Range("H19").Value = 77647

# ---- Row 20: سود (زیان) خالص عملیات در حال تداوم قبل از مالیات ----
$ws.Range("D20").Value = 158639
$ws.Range("E20").Value = 125156
$ws.Range("F20").Value = 270566
$ws.Range("G20").Value = 94524
$ws.Range("H20").Value = 92116

# ---- Row 21: مالیات (Tax) ----
$ws.Range("D21").Value = -6816
$ws.Range("E21").Value = -18
$ws.Range("F21").Value = "-"
$ws.Range("G21").Value = -16267
$ws.Range("H21").Value = -2119

# ---- Row 22: سود (زیان) خالص عملیات در حال تداوم ----
$ws.Range("D22").Value = 151824
$ws.Range("E22").Value = 125137
$ws.Range("F22").Value = 270566
$ws.Range("G22").Value = 78257
$ws.Range("H22").Value = 89997

# Row 23 (سود (زیان) عملیات متوقف شده پس از اثر مالیاتی) stays "-" - no change

# ---- Row 24: سود (زیان) خالص (Net profit) ----
$ws.Range("D24").Value = 151824
$ws.Range("E24").Value = 125137
$ws.Range("F24").Value = 270566
$ws.Range("G24").Value = 78257
$ws.Range("H24").Value = 89997

# ---- Row 25: سود هر سهم پس از کسر مالیات ----
$ws.Range("D25").Value = "-"
$ws.Range("E25").Value = 0
$ws.Range("F25").Value = 0
$ws.Range("G25").Value = 0
$ws.Range("H25").Value = 0

# ---- Row 26: سرمایه (Capital) ----
$ws.Range("D26").Value = 9391
$ws.Range("E26").Value = 7406
$ws.Range("F26").Value = 4202
$ws.Range("G26").Value = 3600
$ws.Range("H26").Value = 87704

# Row 27 (سود هر سهم بر اساس آخرین سرمایه) stays 0 across the board - no change
